$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / add header labels (row 1) ---
$ws.Range("A1").Value = "NOMBRE COMPLETO"
$ws.Range("B1").Value = "DPI / PASAPORTE"
$ws.Range("D1").Value = "GENERO"
$ws.Range("E1").Value = "PAIS DE ORIGEN"
$ws.Range("C1").Value = "EDAD"

# --- Apply the existing blue header style (copied from B1) to the new header cells ---
$ws.Range("B1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("C1").NumberFormat = "0.00"    # EDAD column header uses a numeric format, not text

# --- Give column C (EDAD) a numeric default format / column width ---
$ws.Columns.Item(3).NumberFormat = "0.00"
$ws.Columns.Item(3).ColumnWidth = 11.5546875

# --- Column widths for the other columns ---
$ws.Columns.Item(1).ColumnWidth = 27.88671875
$ws.Columns.Item(2).ColumnWidth = 22.6640625
$ws.Columns.Item(5).ColumnWidth = 14.21875

$ws.Range("D1").Select()
